# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.355.08"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "3.406.69"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "561.27"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").Value = "175.99"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("D8").Value = "3.396.99"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +13.33%  "
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").Value = "54.95"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").Value = "0.0000281"
$ws.Range("E13").Value = "  +6.50%  "
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("D15").Value = "3.950.27"
$ws.Range("D16").Value = "18.40"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.119"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.400.42"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "65.386.41"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "11.90"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").Value = "472.42"
$ws.Range("E22").Value = "  +15.79%  "
$ws.Range("D23").Value = "5.02"
$ws.Range("E23").Value = "  +15.85%  "
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").Value = "87.32"
$ws.Range("E25").Value = "  +5.46%  "
$ws.Range("D26").Value = "13.48"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.92"
$ws.Range("E27").Value = "  +7.62%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.94"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("D29").Value = "8.86"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").Value = "31.35"
$ws.Range("E30").Value = "  +8.43%  "
$ws.Range("D31").Value = "6.76"
$ws.Range("E31").Value = "  +6.86%  "
$ws.Range("D32").Value = "11.57"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").Value = "62.26"
$ws.Range("E33").Value = "  +8.11%  "
$ws.Range("D34").Value = "576.29"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "3.55"
$ws.Range("E37").Value = "  +4.47%  "
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0768"
$ws.Range("E39").Value = "  +4.77%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "35.96"
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").Value = "3.104.05"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("D45").Value = "0.0418"
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").Value = "2.49"
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("D47").Value = "0.136"
$ws.Range("E47").Value = "  +6.57%  "
$ws.Range("D48").Value = "3.16"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "137.55"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("D51").Value = "8.28"
$ws.Range("E51").Value = "  +3.10%  "
